{"js": "// Update the date heading paragraph.\nconst paras = context.document.body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\nparas.items[0].getRange().insertText(\"2024-04-01 Monday\", \"Replace\");\n\n// Update the worksheet table cells (5 blocks of 5 division problems each,\n// separated by blank rows). Replacements are positional (row, col) so the\n// duplicate \"84\u00f77=\" source values resolve to their correct distinct targets.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst newValues = {\n  0: [\"42\u00f78=\", \"70\u00f79=\", \"42\u00f76=\", \"62\u00f73=\", \"65\u00f72=\"],\n  4: [\"44\u00f75=\", \"58\u00f72=\", \"89\u00f77=\", \"76\u00f72=\", \"18\u00f79=\"],\n  8: [\"82\u00f76=\", \"86\u00f73=\", \"14\u00f72=\", \"89\u00f77=\", \"81\u00f79=\"],\n  12: [\"66\u00f73=\", \"60\u00f73=\", \"44\u00f77=\", \"46\u00f73=\", \"83\u00f77=\"],\n  16: [\"61\u00f73=\", \"58\u00f79=\", \"90\u00f73=\", \"22\u00f74=\", \"28\u00f73=\"],\n};\n\nfor (const rowStr of Object.keys(newValues)) {\n  const row = Number(rowStr);\n  const vals = newValues[rowStr];\n  for (let col = 0; col < vals.length; col++) {\n    table.getCell(row, col).value = vals[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date heading paragraph.\n$d.Paragraphs(1).Range.Text = \"2024-04-01 Monday\"\n\n# Update the worksheet table cells (5 blocks of 5 division problems each,\n# separated by blank rows). Replacements are positional (row, col) so the\n# duplicate \"84\u00f77=\" source values resolve to their correct distinct targets.\n$t = $d.Tables(1)\n\n$newValues = @{\n    1  = @(\"42\u00f78=\", \"70\u00f79=\", \"42\u00f76=\", \"62\u00f73=\", \"65\u00f72=\")\n    5  = @(\"44\u00f75=\", \"58\u00f72=\", \"89\u00f77=\", \"76\u00f72=\", \"18\u00f79=\")\n    9  = @(\"82\u00f76=\", \"86\u00f73=\", \"14\u00f72=\", \"89\u00f77=\", \"81\u00f79=\")\n    13 = @(\"66\u00f73=\", \"60\u00f73=\", \"44\u00f77=\", \"46\u00f73=\", \"83\u00f77=\")\n    17 = @(\"61\u00f73=\", \"58\u00f79=\", \"90\u00f73=\", \"22\u00f74=\", \"28\u00f73=\")\n}\n\nforeach ($row in @(1, 5, 9, 13, 17)) {\n    $vals = $newValues[$row]\n    for ($col = 1; $col -le 5; $col++) {\n        $t.Cell($row, $col).Range.Text = $vals[$col - 1]\n    }\n}\n"}
